$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44203
$ws.Cells.Item(2, 11).Value = 25000
$ws.Cells.Item(2, 12).Value = 26000
$ws.Cells.Item(2, 13).Value = 25500
$ws.Cells.Item(2, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 16).Value = 1020

# Row 3
$ws.Cells.Item(3, 4).Value = 44251
$ws.Cells.Item(3, 11).Value = 27000
$ws.Cells.Item(3, 12).Value = 28000
$ws.Cells.Item(3, 13).Value = 27500
$ws.Cells.Item(3, 16).Value = 1100

# Row 4
$ws.Cells.Item(4, 4).Value = 44279
$ws.Cells.Item(4, 11).Value = 28000
$ws.Cells.Item(4, 12).Value = 30000
$ws.Cells.Item(4, 13).Value = 29000
$ws.Cells.Item(4, 15).Value = 'Región del Maule'
$ws.Cells.Item(4, 16).Value = 1160

# Row 5
$ws.Cells.Item(5, 4).Value = 44294
$ws.Cells.Item(5, 11).Value = 20000
$ws.Cells.Item(5, 12).Value = 22000
$ws.Cells.Item(5, 13).Value = 21000
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 840

# Row 6
$ws.Cells.Item(6, 4).Value = 44320
$ws.Cells.Item(6, 11).Value = 26000
$ws.Cells.Item(6, 12).Value = 28000
$ws.Cells.Item(6, 13).Value = 27000
$ws.Cells.Item(6, 16).Value = 1080

# Row 7
$ws.Cells.Item(7, 4).Value = 44194
$ws.Cells.Item(7, 11).Value = 30000
$ws.Cells.Item(7, 12).Value = 32000
$ws.Cells.Item(7, 13).Value = 31000
$ws.Cells.Item(7, 16).Value = 1240

# Row 8
$ws.Cells.Item(8, 4).Value = 44349
$ws.Cells.Item(8, 10).Value = 60
$ws.Cells.Item(8, 11).Value = 30000
$ws.Cells.Item(8, 12).Value = 32000
$ws.Cells.Item(8, 13).Value = 31000
$ws.Cells.Item(8, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(8, 16).Value = 1240

# Row 9
$ws.Cells.Item(9, 4).Value = 44265
$ws.Cells.Item(9, 11).Value = 22000
$ws.Cells.Item(9, 12).Value = 24000
$ws.Cells.Item(9, 13).Value = 23000
$ws.Cells.Item(9, 16).Value = 920

# Row 10
$ws.Cells.Item(10, 4).Value = 44328
$ws.Cells.Item(10, 11).Value = 32000
$ws.Cells.Item(10, 12).Value = 34000
$ws.Cells.Item(10, 13).Value = 33000
$ws.Cells.Item(10, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(10, 16).Value = 1320

# Row 11
$ws.Cells.Item(11, 4).Value = 44216
$ws.Cells.Item(11, 11).Value = 26000
$ws.Cells.Item(11, 12).Value = 28000
$ws.Cells.Item(11, 13).Value = 27000
$ws.Cells.Item(11, 16).Value = 1080

# Row 12
$ws.Cells.Item(12, 4).Value = 44272
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 22000
$ws.Cells.Item(12, 12).Value = 24000
$ws.Cells.Item(12, 13).Value = 23000
$ws.Cells.Item(12, 15).Value = 'Región del Maule'
$ws.Cells.Item(12, 16).Value = 920

# Row 13
$ws.Cells.Item(13, 4).Value = 44308
$ws.Cells.Item(13, 11).Value = 28000
$ws.Cells.Item(13, 12).Value = 30000
$ws.Cells.Item(13, 13).Value = 29000
$ws.Cells.Item(13, 16).Value = 1160

# Row 14
$ws.Cells.Item(14, 4).Value = 44236
$ws.Cells.Item(14, 11).Value = 25000
$ws.Cells.Item(14, 12).Value = 26000
$ws.Cells.Item(14, 13).Value = 25500
$ws.Cells.Item(14, 15).Value = 'Región del Maule'
$ws.Cells.Item(14, 16).Value = 1020

# Row 15
$ws.Cells.Item(15, 4).Value = 44188
$ws.Cells.Item(15, 11).Value = 42000
$ws.Cells.Item(15, 12).Value = 44000
$ws.Cells.Item(15, 13).Value = 43000
$ws.Cells.Item(15, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(15, 16).Value = 1720

# Row 17
$ws.Cells.Item(17, 4).Value = 44313
$ws.Cells.Item(17, 11).Value = 30000
$ws.Cells.Item(17, 12).Value = 32000
$ws.Cells.Item(17, 13).Value = 31000
$ws.Cells.Item(17, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 16).Value = 1240

# Row 18
$ws.Cells.Item(18, 4).Value = 44244
$ws.Cells.Item(18, 11).Value = 25000
$ws.Cells.Item(18, 12).Value = 26000
$ws.Cells.Item(18, 13).Value = 25500
$ws.Cells.Item(18, 16).Value = 1020

# Row 19
$ws.Cells.Item(19, 4).Value = 44316
$ws.Cells.Item(19, 11).Value = 26000
$ws.Cells.Item(19, 12).Value = 27000
$ws.Cells.Item(19, 13).Value = 26500
$ws.Cells.Item(19, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 16).Value = 1060

# Row 20
$ws.Cells.Item(20, 4).Value = 44230
$ws.Cells.Item(20, 11).Value = 35000
$ws.Cells.Item(20, 12).Value = 36000
$ws.Cells.Item(20, 13).Value = 35500
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 1420

# Row 21
$ws.Cells.Item(21, 4).Value = 44342
$ws.Cells.Item(21, 11).Value = 28000
$ws.Cells.Item(21, 12).Value = 30000
$ws.Cells.Item(21, 13).Value = 29000
$ws.Cells.Item(21, 16).Value = 1160
